# Apply changes described in the commit "fixing merge problems in post and main"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the search text value in P2 from "Ja" to "Jane"
$ws.Range("P2").Value = "Jane"

# Update the visible window / selection state of the sheet view
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 12   # column L is left-most visible column
$ws.Range("P2").Select()

$wb.Save()
